# Padding-Margin.docx edit script
#
# 1. Fix typo: "rectangular boxes is described" -> "rectangular box is described"
# 2. Remove the stray "_GoBack" bookmark that wraps nothing (right after "Content ... ,")
#    (this also naturally renumbers the surviving "_Hlk15334437" bookmark id from 1 to 0,
#    exactly like Word does when a lower-numbered bookmark id is freed)
# 3. Reword the closing sentence of the Margin section ("...on left and on right." ->
#    "...on the left and on the right.") splitting it into several runs the way Word's
#    editor naturally does, and drop a fresh "_GoBack" bookmark (Word always leaves one
#    at the site of the most recent edit) at the point right before " right."

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "rectangular boxes is described" -> "rectangular box is described"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "each of these rectangular boxes is described", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "each of these rectangular box is described", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Drop the old "_GoBack" bookmark (the pair right after "...Element Box,")
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 3) Reword "...on left and on right." -> "...on the left and on the right."
#    and split the sentence into the same run layout as the real edit, with a
#    "_GoBack" bookmark dropped right before the final " right." run.
# ---------------------------------------------------------------------------
$old = "The margin above will provide 100px of spacing outside the border from four sides, which mean the box will have 100px of empty space on the top, on the bottom, on left and on right."
$new = "The margin above will provide 100px of spacing outside the border from four sides, which mean the box will have 100px of empty space on the top, on the bottom, on the left and on the right."

$found = $d.Content
$found.Find.Execute($old) | Out-Null
$sentenceStart = $found.Start

$target = $d.Range($sentenceStart, $found.End)
$target.Text = $new

# Offsets (relative to $sentenceStart) of the run boundaries in the new text:
#   "...on the bottom, on" | " the" | " left and on" | " the" | <bookmark> | " right."
$off1 = 162   # end of "...on"              / start of " the"
$off2 = 166   # end of " the"               / start of " left and on"
$off3 = 178   # end of " left and on"       / start of " the"
$off4 = 182   # end of " the"               / bookmark / start of " right."

# Force run splits at off1/off2/off3 using disposable bookmarks: adding a
# bookmark at a collapsed point splits the run there, and the split survives
# even after the temporary bookmark is deleted again.
$d.Bookmarks.Add("_zz_tmp1", $d.Range($sentenceStart + $off1, $sentenceStart + $off1)) | Out-Null
$d.Bookmarks.Add("_zz_tmp2", $d.Range($sentenceStart + $off2, $sentenceStart + $off2)) | Out-Null
$d.Bookmarks.Add("_zz_tmp3", $d.Range($sentenceStart + $off3, $sentenceStart + $off3)) | Out-Null

# The real bookmark Word leaves behind at the edit point.
$d.Bookmarks.Add("_GoBack", $d.Range($sentenceStart + $off4, $sentenceStart + $off4)) | Out-Null

$d.Bookmarks("_zz_tmp1").Delete()
$d.Bookmarks("_zz_tmp2").Delete()
$d.Bookmarks("_zz_tmp3").Delete()
